{"js": "const body = context.document.body;\nconst doc = context.document;\n\n// ---------------------------------------------------------------------\n// 1) Title: \"...running ETL(Kettle) Job\" -> \"...running ETL (Kettle) Job\"\n//    (a space is added between \"ETL\" and \"(Kettle)\"). Replacing the whole\n//    paragraph range (not just the inner run) cleanly drops the now-stale\n//    gramStart/gramEnd proofing markers that used to flag the missing\n//    space, instead of leaving them orphaned in the middle of the text.\n// ---------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items.find(\n  (p) => p.text.indexOf(\"running ETL(Kettle) Job\") !== -1\n);\n\nif (titleParagraph) {\n  const titleRange = titleParagraph.getRange(\"Whole\");\n  titleRange.insertText(\n    \"Creating new database and running ETL (Kettle) Job\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) \"Copy the ETL folder from\" -> \"Copy the source/database/ETL folder from\"\n// ---------------------------------------------------------------------\nconst copyResults = body.search(\"Copy the ETL folder from\", { matchCase: true });\ncopyResults.load(\"items\");\nawait context.sync();\nif (copyResults.items.length > 0) {\n  copyResults.items[0].insertText(\n    \"Copy the source/database/ETL folder from\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) \"If the ETL folder from\" -> \"If the source/database/ETL folder from\"\n// ---------------------------------------------------------------------\nconst ifResults = body.search(\"If the ETL folder from\", { matchCase: true });\nifResults.load(\"items\");\nawait context.sync();\nif (ifResults.items.length > 0) {\n  ifResults.items[0].insertText(\n    \"If the source/database/ETL folder from\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 4) Move the \"_GoBack\" bookmark: it used to sit inside \"Java Runtime\"\n//    (splitting \"Java Runtim\" / \"e\"); it now belongs in the title, right\n//    between \"ETL \" and \"(Kettle)\" - i.e. at the point we just edited.\n// ---------------------------------------------------------------------\nconst existingBookmark = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\nexistingBookmark.load(\"isNullObject\");\nawait context.sync();\nif (!existingBookmark.isNullObject) {\n  doc.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\nconst retitled = body.paragraphs;\nretitled.load(\"items/text\");\nawait context.sync();\nconst newTitleParagraph = retitled.items.find(\n  (p) => p.text.indexOf(\"running ETL (Kettle) Job\") !== -1\n);\nif (newTitleParagraph) {\n  const kettleResults = newTitleParagraph.search(\"(Kettle)\", { matchCase: true });\n  kettleResults.load(\"items\");\n  await context.sync();\n  if (kettleResults.items.length > 0) {\n    const insertionPoint = kettleResults.items[0].getRange(Word.RangeLocation.start);\n    insertionPoint.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Title: \"...running ETL(Kettle) Job\" -> \"...running ETL (Kettle) Job\"\n#    (a space is added between \"ETL\" and \"(Kettle)\"). Replacing the whole\n#    paragraph's text span (start..end, i.e. not including the paragraph\n#    mark) cleanly drops the now-stale gramStart/gramEnd proofing markers\n#    that used to flag the missing space, instead of leaving them\n#    orphaned in the middle of the text.\n# ---------------------------------------------------------------------\n$titlePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*running ETL(Kettle) Job*\") {\n        $titlePara = $p\n        break\n    }\n}\n\nif ($titlePara -ne $null) {\n    $pStart = $titlePara.Range.Start\n    $pEnd = $titlePara.Range.End\n    $titleRange = $d.Range($pStart, $pEnd)\n    $titleRange.Text = \"Creating new database and running ETL (Kettle) Job\"\n}\n\n# ---------------------------------------------------------------------\n# 2) \"Copy the ETL folder from\" -> \"Copy the source/database/ETL folder from\"\n# ---------------------------------------------------------------------\n$find1 = $d.Content\n$found1 = $find1.Find.Execute(\"Copy the ETL folder from\")\nif ($found1) {\n    $r1 = $d.Range($find1.Start, $find1.End)\n    $r1.Text = \"Copy the source/database/ETL folder from\"\n}\n\n# ---------------------------------------------------------------------\n# 3) \"If the ETL folder from\" -> \"If the source/database/ETL folder from\"\n# ---------------------------------------------------------------------\n$find2 = $d.Content\n$found2 = $find2.Find.Execute(\"If the ETL folder from\")\nif ($found2) {\n    $r2 = $d.Range($find2.Start, $find2.End)\n    $r2.Text = \"If the source/database/ETL folder from\"\n}\n\n# ---------------------------------------------------------------------\n# 4) Move the \"_GoBack\" bookmark: it used to sit inside \"Java Runtime\"\n#    (splitting \"Java Runtim\" / \"e\"); it now belongs in the title, right\n#    between \"ETL \" and \"(Kettle)\" - i.e. at the point we just edited.\n# ---------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$find3 = $d.Content\n$found3 = $find3.Find.Execute(\"(Kettle)\")\nif ($found3) {\n    $insertPoint = $d.Range($find3.Start, $find3.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $insertPoint)\n}\n"}
